# Add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (before the current
#    "2022-Q3" sheet), cloned from "2022-Q3" so that all cell styles
#    (header border/bold, index-column style, etc.) come along for free.
# 2. Trim it down to the 4 data rows for 2022-Q4 and overwrite the values.
# 3. Update the "总计" summary sheet: shift the existing quarter rows down
#    by one and insert the new 2022-Q4 summary row at the top.
# 4. Restore the originally-active sheet (last tab) as the active sheet.

$wb = $excel.ActiveWorkbook

# Helper: write a value as literal text (never let Excel's COM layer
# re-interpret a numeric-looking string like "6.60" as a number), while
# leaving the cell's style untouched (reset back to Normal afterwards so we
# don't leave a stray "@" number format behind).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1+2: create the "2022-Q4" sheet as a clone of "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)                      # clone lands right before the original
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# 2022-Q3 has 5 data rows (rows 2-6); 2022-Q4 only needs 4 (rows 2-5).
$q4.Rows.Item(6).Delete()

Set-TextValue $q4.Range("B2") "014202"
Set-TextValue $q4.Range("C2") "天弘中证1000指数增强C"
Set-TextValue $q4.Range("D2") "6.60"
Set-TextValue $q4.Range("E2") "94.11"
Set-TextValue $q4.Range("F2") "1.47"
Set-TextValue $q4.Range("G2") "0.0970"
$q4.Range("H2").Value = 10

Set-TextValue $q4.Range("B3") "014201"
Set-TextValue $q4.Range("C3") "天弘中证1000指数增强A"
Set-TextValue $q4.Range("D3") "3.86"
Set-TextValue $q4.Range("E3") "94.11"
Set-TextValue $q4.Range("F3") "1.47"
Set-TextValue $q4.Range("G3") "0.0567"
$q4.Range("H3").Value = 10

Set-TextValue $q4.Range("B4") "000926"
Set-TextValue $q4.Range("C4") "中信建投睿信灵活配置混合A"
Set-TextValue $q4.Range("D4") "0.10"
Set-TextValue $q4.Range("E4") "83.25"
Set-TextValue $q4.Range("F4") "1.07"
Set-TextValue $q4.Range("G4") "0.0011"
$q4.Range("H4").Value = 3

Set-TextValue $q4.Range("B5") "004676"
Set-TextValue $q4.Range("C5") "中信建投睿信灵活配置混合C"
Set-TextValue $q4.Range("D5") "0.03"
Set-TextValue $q4.Range("E5") "83.25"
Set-TextValue $q4.Range("F5") "1.07"
Set-TextValue $q4.Range("G5") "0.0003"
$q4.Range("H5").Value = 3

# ---------------------------------------------------------------------
# Step 3: update "总计" - push the existing rows down one and write the
# new 2022-Q4 summary row on top. Column A is just the 0-based row index
# (row - 2) and is already correct in every existing row, so only B:D are
# shifted; the new row 8 needs a fresh A8 cell (copied for its style).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A7").Copy($total.Range("A8"))
$total.Range("A8").Value = 6

$total.Range("B7:D7").Copy($total.Range("B8:D8"))
$total.Range("B6:D6").Copy($total.Range("B7:D7"))
$total.Range("B5:D5").Copy($total.Range("B6:D6"))
$total.Range("B4:D4").Copy($total.Range("B5:D5"))
$total.Range("B3:D3").Copy($total.Range("B4:D4"))
$total.Range("B2:D2").Copy($total.Range("B3:D3"))

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.16

# ---------------------------------------------------------------------
# Step 4: restore the active sheet to the last tab (matches original file)
# ---------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
